$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = "@"
$c.Value = '62.802.27'
$c.Style = "Normal"
$c = $ws.Range('E2')
$c.NumberFormat = "@"
$c.Value = '  -0.52%  '
$c.Style = "Normal"

$c = $ws.Range('D3')
$c.NumberFormat = "@"
$c.Value = '2.580.27'
$c.Style = "Normal"
$c = $ws.Range('E3')
$c.NumberFormat = "@"
$c.Value = '  +1.11%  '
$c.Style = "Normal"

$c = $ws.Range('E4')
$c.NumberFormat = "@"
$c.Value = '  +0.02%  '
$c.Style = "Normal"

$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '583.51'
$c.Style = "Normal"
$c = $ws.Range('E5')
$c.NumberFormat = "@"
$c.Value = '  -0.24%  '
$c.Style = "Normal"

$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '145.96'
$c.Style = "Normal"
$c = $ws.Range('E6')
$c.NumberFormat = "@"
$c.Value = '  -1.10%  '
$c.Style = "Normal"

$c = $ws.Range('E7')
$c.NumberFormat = "@"
$c.Value = '  +0.01%  '
$c.Style = "Normal"

$c = $ws.Range('E8')
$c.NumberFormat = "@"
$c.Value = '  +1.47%  '
$c.Style = "Normal"

$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '0.107'
$c.Style = "Normal"
$c = $ws.Range('E9')
$c.NumberFormat = "@"
$c.Value = '  +1.12%  '
$c.Style = "Normal"

$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '5.62'
$c.Style = "Normal"
$c = $ws.Range('E10')
$c.NumberFormat = "@"
$c.Value = '  +1.35%  '
$c.Style = "Normal"

$c = $ws.Range('E11')
$c.NumberFormat = "@"
$c.Value = '  -0.18%  '
$c.Style = "Normal"

$c = $ws.Range('E12')
$c.NumberFormat = "@"
$c.Value = '  -0.82%  '
$c.Style = "Normal"

$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value = '27.13'
$c.Style = "Normal"
$c = $ws.Range('E13')
$c.NumberFormat = "@"
$c.Value = '  -1.25%  '
$c.Style = "Normal"

$c = $ws.Range('D14')
$c.NumberFormat = "@"
$c.Value = '3.042.74'
$c.Style = "Normal"
$c = $ws.Range('E14')
$c.NumberFormat = "@"
$c.Value = '  +1.30%  '
$c.Style = "Normal"

$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '62.706.82'
$c.Style = "Normal"
$c = $ws.Range('E15')
$c.NumberFormat = "@"
$c.Value = '  -0.49%  '
$c.Style = "Normal"

$c = $ws.Range('E16')
$c.NumberFormat = "@"
$c.Value = '  +1.51%  '
$c.Style = "Normal"

$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '2.579.10'
$c.Style = "Normal"
$c = $ws.Range('E17')
$c.NumberFormat = "@"
$c.Value = '  +0.64%  '
$c.Style = "Normal"

$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '11.25'
$c.Style = "Normal"
$c = $ws.Range('E18')
$c.NumberFormat = "@"
$c.Value = '  -1.10%  '
$c.Style = "Normal"

$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '340.61'
$c.Style = "Normal"
$c = $ws.Range('E19')
$c.NumberFormat = "@"
$c.Value = '  +0.97%  '
$c.Style = "Normal"

$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '4.37'
$c.Style = "Normal"
$c = $ws.Range('E20')
$c.NumberFormat = "@"
$c.Value = '  +0.84%  '
$c.Style = "Normal"

$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '6.67'
$c.Style = "Normal"
$c = $ws.Range('E21')
$c.NumberFormat = "@"
$c.Value = '  -1.76%  '
$c.Style = "Normal"

$c = $ws.Range('E22')
$c.NumberFormat = "@"
$c.Value = '  -0.32%  '
$c.Style = "Normal"

$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '67.31'
$c.Style = "Normal"
$c = $ws.Range('E23')
$c.NumberFormat = "@"
$c.Value = '  +2.10%  '
$c.Style = "Normal"

$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '2.715.44'
$c.Style = "Normal"
$c = $ws.Range('E24')
$c.NumberFormat = "@"
$c.Value = '  +1.46%  '
$c.Style = "Normal"

$c = $ws.Range('E25')
$c.NumberFormat = "@"
$c.Value = '  -2.23%  '
$c.Style = "Normal"

$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '1.58'
$c.Style = "Normal"
$c = $ws.Range('E26')
$c.NumberFormat = "@"
$c.Value = '  -2.49%  '
$c.Style = "Normal"

$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$c = $ws.Range('E27')
$c.NumberFormat = "@"
$c.Value = '  +0.14%  '
$c.Style = "Normal"

$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '7.87'
$c.Style = "Normal"
$c = $ws.Range('E28')
$c.NumberFormat = "@"
$c.Value = '  +2.35%  '
$c.Style = "Normal"

$c = $ws.Range('E29')
$c.NumberFormat = "@"
$c.Value = '  -1.89%  '
$c.Style = "Normal"

$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '8.27'
$c.Style = "Normal"
$c = $ws.Range('E30')
$c.NumberFormat = "@"
$c.Value = '  -1.82%  '
$c.Style = "Normal"

$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '1.92'
$c.Style = "Normal"
$c = $ws.Range('E31')
$c.NumberFormat = "@"
$c.Value = '  -2.57%  '
$c.Style = "Normal"

$c = $ws.Range('B32')
$c.NumberFormat = "@"
$c.Value = 'Bittensor'
$c.Style = "Normal"
$c = $ws.Range('C32')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$c.Style = "Normal"
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '466.69'
$c.Style = "Normal"
$c = $ws.Range('E32')
$c.NumberFormat = "@"
$c.Value = '  +11.21%  '
$c.Style = "Normal"

$c = $ws.Range('B33')
$c.NumberFormat = "@"
$c.Value = 'PEPE'
$c.Style = "Normal"
$c = $ws.Range('C33')
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$c.Style = "Normal"
$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '0.0₃0816'
$c.Style = "Normal"
$c = $ws.Range('E33')
$c.NumberFormat = "@"
$c.Value = '  -0.26%  '
$c.Style = "Normal"

$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '176.41'
$c.Style = "Normal"
$c = $ws.Range('E34')
$c.NumberFormat = "@"
$c.Value = '  -1.06%  '
$c.Style = "Normal"

$c = $ws.Range('E35')
$c.NumberFormat = "@"
$c.Value = '  +2.85%  '
$c.Style = "Normal"

$c = $ws.Range('E36')
$c.NumberFormat = "@"
$c.Value = '  +0.05%  '
$c.Style = "Normal"

$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '0.399'
$c.Style = "Normal"
$c = $ws.Range('E37')
$c.NumberFormat = "@"
$c.Value = '  -0.67%  '
$c.Style = "Normal"

$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '18.95'
$c.Style = "Normal"
$c = $ws.Range('E38')
$c.NumberFormat = "@"
$c.Value = '  -1.29%  '
$c.Style = "Normal"

$c = $ws.Range('E39')
$c.NumberFormat = "@"
$c.Value = '  +2.60%  '
$c.Style = "Normal"

$c = $ws.Range('E40')
$c.NumberFormat = "@"
$c.Value = '  +0.00%  '
$c.Style = "Normal"

$c = $ws.Range('E41')
$c.NumberFormat = "@"
$c.Value = '  -3.30%  '
$c.Style = "Normal"

$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '157.98'
$c.Style = "Normal"
$c = $ws.Range('E42')
$c.NumberFormat = "@"
$c.Value = '  +4.89%  '
$c.Style = "Normal"

$c = $ws.Range('E43')
$c.NumberFormat = "@"
$c.Value = '  -1.90%  '
$c.Style = "Normal"

$c = $ws.Range('E44')
$c.NumberFormat = "@"
$c.Value = '  +4.10%  '
$c.Style = "Normal"

$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '20.93'
$c.Style = "Normal"
$c = $ws.Range('E45')
$c.NumberFormat = "@"
$c.Value = '  +0.16%  '
$c.Style = "Normal"

$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '0.0538'
$c.Style = "Normal"
$c = $ws.Range('E46')
$c.NumberFormat = "@"
$c.Value = '  -1.05%  '
$c.Style = "Normal"

$c = $ws.Range('E47')
$c.NumberFormat = "@"
$c.Value = '  -0.87%  '
$c.Style = "Normal"

$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '0.0235'
$c.Style = "Normal"
$c = $ws.Range('E48')
$c.NumberFormat = "@"
$c.Value = '  -1.35%  '
$c.Style = "Normal"

$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '18.24'
$c.Style = "Normal"
$c = $ws.Range('E49')
$c.NumberFormat = "@"
$c.Value = '  -0.46%  '
$c.Style = "Normal"

$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '11.41'
$c.Style = "Normal"
$c = $ws.Range('E50')
$c.NumberFormat = "@"
$c.Value = '  +0.96%  '
$c.Style = "Normal"

$c = $ws.Range('E51')
$c.NumberFormat = "@"
$c.Value = '  -0.86%  '
$c.Style = "Normal"
